$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metadata record (row 2) for MCH146
$ws.Rows.Item(2).RowHeight = 15.75
$ws.Range("A2").Value = "MCH146"
$ws.Range("C2").Value = "NAMIBIA FACTS AND FIGURES, THE CONSTITUTION, NATIONAL FLAG, STATISTICAL/ECONOMIC REVIEW"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# Match the font used for the rest of the data rows (Calibri 10pt, theme text color).
# D2/H2 pick up the same formatting even though they stay empty, like the source file.
$dataRange = $ws.Range("A2:H2")
$dataRange.Font.Name = "Calibri"
$dataRange.Font.ThemeColor = 1

# There is no alternativeIdentifiers value for this record - leave B2 untouched/empty
$ws.Range("B2").Clear()

# Keep the header row frozen and select the new data row, like the source file
$ws.Range("A2:I2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
